# Femacal de La Calera - Repollo: insert two new weekly price records
# (date 44516, qualities "Primera" and "Segunda") right before the
# existing row for date 44270, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 379-380; everything from the old row 379
# onward (through the old row 399) shifts down to 381-401, and the
# sheet's dimension/used-range grows from R399 to R401 automatically.
$ws.Rows("379:380").Insert()

# New row 379: Primera, 2021-11-16 (serial 44516)
$ws.Cells.Item(379, 1).Value  = 3
$ws.Cells.Item(379, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(379, 3).Value  = "Coquimbo"
$ws.Cells.Item(379, 4).Value  = 44516
$ws.Cells.Item(379, 5).Value  = 5
$ws.Cells.Item(379, 6).Value  = 100112006
$ws.Cells.Item(379, 7).Value  = "Repollo"
$ws.Cells.Item(379, 8).Value  = "Crespo record"
$ws.Cells.Item(379, 9).Value  = "Primera"
$ws.Cells.Item(379, 10).Value = 2600
$ws.Cells.Item(379, 11).Value = 600
$ws.Cells.Item(379, 12).Value = 650
$ws.Cells.Item(379, 13).Value = 625
$ws.Cells.Item(379, 14).Value = "$/unidad"
$ws.Cells.Item(379, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(379, 16).Value = 625
$ws.Cells.Item(379, 17).Value = 1
$ws.Cells.Item(379, 18).Value = "Hortaliza"

# New row 380: Segunda, same date (serial 44516)
$ws.Cells.Item(380, 1).Value  = 3
$ws.Cells.Item(380, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(380, 3).Value  = "Coquimbo"
$ws.Cells.Item(380, 4).Value  = 44516
$ws.Cells.Item(380, 5).Value  = 5
$ws.Cells.Item(380, 6).Value  = 100112006
$ws.Cells.Item(380, 7).Value  = "Repollo"
$ws.Cells.Item(380, 8).Value  = "Crespo record"
$ws.Cells.Item(380, 9).Value  = "Segunda"
$ws.Cells.Item(380, 10).Value = 1100
$ws.Cells.Item(380, 11).Value = 500
$ws.Cells.Item(380, 12).Value = 500
$ws.Cells.Item(380, 13).Value = 500
$ws.Cells.Item(380, 14).Value = "$/unidad"
$ws.Cells.Item(380, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(380, 16).Value = 500
$ws.Cells.Item(380, 17).Value = 1
$ws.Cells.Item(380, 18).Value = "Hortaliza"
